$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new recognizer row
$ws.Range("A6").Value = "auto_receipt"
$ws.Range("B6").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog"

# Turn B6 into a hyperlink (mirrors the existing rows' hyperlinks)
$ws.Hyperlinks.Add($ws.Range("B6"), "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog") | Out-Null

# Match styling used by the other hyperlink cells (B2:B5)
$ws.Range("B6").Style = $ws.Range("B5").Style

# Column B widens (bestFit) to accommodate the newly added, longer URL text
$ws.Columns.Item(2).ColumnWidth = 53.06640625

# Reflect the new selection shown in the saved workbook
$ws.Range("C6:D6").Select() | Out-Null
